# Updated cryptos list on Wed Jul 24 04:58:13 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to hold a literal text value (not auto-coerced to a
# number) while leaving the cell's style/formatting exactly as it was
# before (the source sheet keeps every Price/Volume cell as plain,
# unstyled text, even when the text looks numeric, e.g. "583.12").
function Set-TextValue($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws "D2" "65.850.51"
Set-TextValue $ws "E2" "  -1.85%  "

# Row 3 - Ethereum
Set-TextValue $ws "D3" "3.438.80"
Set-TextValue $ws "E3" "  -0.59%  "

# Row 5 - BNB
Set-TextValue $ws "D5" "583.07"
Set-TextValue $ws "E5" "  -0.30%  "

# Row 6 - Solana
Set-TextValue $ws "D6" "173.45"
Set-TextValue $ws "E6" "  -1.52%  "

# Row 7 - USDC
Set-TextValue $ws "E7" "  +0.00%  "

# Row 8 - XRP
Set-TextValue $ws "D8" "0.605"
Set-TextValue $ws "E8" "  -0.50%  "

# Row 9 - LidoStakedEther
Set-TextValue $ws "D9" "3.437.53"
Set-TextValue $ws "E9" "  -0.59%  "

# Row 10 - Dogecoin
Set-TextValue $ws "E10" "  -3.05%  "

# Row 11 - Toncoin
Set-TextValue $ws "E11" "  -0.08%  "

# Row 12 - Cardano
Set-TextValue $ws "E12" "  -3.50%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue $ws "D13" "4.033.52"
Set-TextValue $ws "E13" "  -0.72%  "

# Row 14 - TRON
Set-TextValue $ws "E14" "  +1.43%  "

# Row 15 - Avalanche
Set-TextValue $ws "D15" "28.64"
Set-TextValue $ws "E15" "  -9.36%  "

# Row 16 - WrappedBTC
Set-TextValue $ws "D16" "65.912.03"
Set-TextValue $ws "E16" "  -1.83%  "

# Row 17 - ShibaInu
Set-TextValue $ws "E17" "  -2.02%  "

# Row 18 - WrappedEther
Set-TextValue $ws "D18" "3.438.61"
Set-TextValue $ws "E18" "  -0.54%  "

# Row 19 - Polkadot
Set-TextValue $ws "E19" "  -2.35%  "

# Row 20 - Chainlink
Set-TextValue $ws "D20" "13.82"
Set-TextValue $ws "E20" "  -0.63%  "

# Row 21 - BitcoinCash
Set-TextValue $ws "D21" "368.53"
Set-TextValue $ws "E21" "  -2.89%  "

# Row 22 - Uniswap
Set-TextValue $ws "D22" "7.65"
Set-TextValue $ws "E22" "  -2.07%  "

# Row 23 - Litecoin
Set-TextValue $ws "D23" "72.12"
Set-TextValue $ws "E23" "  +0.80%  "

# Row 24 - Dai
Set-TextValue $ws "E24" "  -0.10%  "

# Row 25 - Polygon
Set-TextValue $ws "E25" "  +0.53%  "

# Row 26 - PEPE
Set-TextValue $ws "E26" "  +0.78%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue $ws "E27" "  -2.34%  "

# Row 28 - Kaspa
Set-TextValue $ws "E28" "  +1.68%  "

# Row 29 - Binance-PegBSC-USD
Set-TextValue $ws "D29" "1.00"
Set-TextValue $ws "E29" "  -0.01%  "

# Row 30 - EthereumClassic
Set-TextValue $ws "D30" "23.64"
Set-TextValue $ws "E30" "  -1.37%  "

# Row 31 - NEARProtocol
Set-TextValue $ws "D31" "5.73"
Set-TextValue $ws "E31" "  -3.85%  "

# Row 32 - PancakeSwap
Set-TextValue $ws "D32" "1.97"
Set-TextValue $ws "E32" "  -2.32%  "

# Row 33 - USDe
Set-TextValue $ws "D33" "1.00"
Set-TextValue $ws "E33" "  +0.00%  "

# Row 34 - Fetch.AI
Set-TextValue $ws "E34" "  -5.82%  "

# Row 35 - Aptos
Set-TextValue $ws "D35" "6.99"
Set-TextValue $ws "E35" "  -2.82%  "

# Row 36 - ImmutableX
Set-TextValue $ws "E36" "  -0.61%  "

# Row 37 - Monero
Set-TextValue $ws "D37" "160.72"
Set-TextValue $ws "E37" "  +0.31%  "

# Row 38 - EnergySwap
Set-TextValue $ws "D38" "28.82"
Set-TextValue $ws "E38" "  +6.27%  "

# Row 39 - Mantle
Set-TextValue $ws "E39" "  +0.11%  "

# Row 40/41 - coins swapped (Stacks <-> dogwifhat), with updated price/volume
Set-TextValue $ws "B40" "dogwifhat"
Set-TextValue $ws "C40" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws "D40" "2.61"
Set-TextValue $ws "E40" "  -1.37%  "

Set-TextValue $ws "B41" "Stacks"
Set-TextValue $ws "C41" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws "D41" "1.76"
Set-TextValue $ws "E41" "  -2.80%  "

# Row 42 - Maker
Set-TextValue $ws "D42" "2.771.84"
Set-TextValue $ws "E42" "  +2.64%  "

# Row 43 - RenderToken
Set-TextValue $ws "E43" "  -2.76%  "

# Row 44 - Filecoin
Set-TextValue $ws "E44" "  -0.48%  "

# Row 45 - Hedera
Set-TextValue $ws "E45" "  -2.46%  "

# Row 46 - OKB
Set-TextValue $ws "D46" "40.18"
Set-TextValue $ws "E46" "  -2.17%  "

# Row 47 - InjectiveProtocol
Set-TextValue $ws "D47" "24.75"
Set-TextValue $ws "E47" "  -2.60%  "

# Row 48 - VeChain
Set-TextValue $ws "E48" "  -1.77%  "

# Row 49 - Bittensor
Set-TextValue $ws "D49" "325.46"
Set-TextValue $ws "E49" "  +0.91%  "

# Row 50 - Stellar
Set-TextValue $ws "E50" "  -1.45%  "

# Row 51 - Cosmos
Set-TextValue $ws "D51" "6.26"
Set-TextValue $ws "E51" "  +0.35%  "
